$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.00639
$ws.Range("E2").Value = 0.5329999999999999
$ws.Range("G2").Value = -0.1877271672724461
$ws.Range("H2").Value = -0.1877271672724461
$ws.Range("I2").Value = 1.824800866135643
$ws.Range("J2").Value = 1.823468894700508
$ws.Range("K2").Value = -76.39
$ws.Range("L2").Value = 2.953754543345449
$ws.Range("M2").Value = 3.21066
$ws.Range("N2").Value = 0.004863973094578012
$ws.Range("O2").Value = -0.04202984683859144
$ws.Range("P2").Value = 3.21066
$ws.Range("Q2").Value = 0.004863973094578012
$ws.Range("R2").Value = -0.04202984683859144
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 7.527
$ws.Range("V2").Value = 0.01140299050129528
$ws.Range("W2").Value = -0.09023778462998103
$ws.Range("X2").Value = 0.0437464840940926
$ws.Range("Y2").Value = -0.1339842687240736
$ws.Range("Z2").Value = -0.01909387168079523
$ws.Range("AA2").Value = -0.00415563108796092
$ws.Range("AB2").Value = 0.0437464840940926
$ws.Range("AC2").Value = -0.04930217260413078
$ws.Range("AD2").Value = 56.81
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 56.81
$ws.Range("AG2").Value = 49.283
$ws.Range("AH2").Value = 0.07924396708048542
$ws.Range("AI2").Value = 0.04910494334045
$ws.Range("AJ2").Value = 0.06947402847303182
$ws.Range("AK2").Value = 0.04287778747380117
$ws.Range("AL2").Value = 2.852
$ws.Range("AM2").Value = 2.59
$ws.Range("AN2").Value = 18.09812042051609
$ws.Range("AO2").Value = -16.54733520336606
$ws.Range("AP2").Value = 15.70022300095572
$ws.Range("AQ2").Value = -18.22123552123552

# Row 3
$ws.Range("B3").Value = "National Investment Trust Ltd (MUSE:NITL.N0000)"
$ws.Range("D3").Value = 0.314
$ws.Range("E3").Value = 0.5329999999999999
$ws.Range("G3").Value = 0.7563739376770539
$ws.Range("H3").Value = 0.7563739376770539
$ws.Range("I3").Value = 0.7762039660056659
$ws.Range("J3").Value = 0.7716713881019832
$ws.Range("K3").Value = 2.72
$ws.Range("L3").Value = 0.7705382436260624
$ws.Range("M3").Value = 0.756
$ws.Range("N3").Value = 0.02918918918918919
$ws.Range("O3").Value = 0.2779411764705882
$ws.Range("P3").Value = 0.756
$ws.Range("Q3").Value = 0.02918918918918919
$ws.Range("R3").Value = 0.2779411764705882
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 1.04
$ws.Range("V3").Value = 0.04015444015444015
$ws.Range("W3").Value = 0.09543859649122807
$ws.Range("X3").Value = 0.0437464840940926
$ws.Range("Y3").Value = 0.05169211239713548
$ws.Range("Z3").Value = 0.1316672883252518
$ws.Range("AA3").Value = 0.1016038791495711
$ws.Range("AB3").Value = 0.0437464840940926
$ws.Range("AC3").Value = 0.05785739505547847
$ws.Range("AG3").Value = -1.04
$ws.Range("AJ3").Value = -0.0418342719227675
$ws.Range("AK3").Value = -0.03843311160384331
$ws.Range("AL3").Value = 0.004
$ws.Range("AM3").Value = 0.004
$ws.Range("AN3").Value = 0
$ws.Range("AO3").Value = 685
$ws.Range("AP3").Value = -0.3661971830985916
$ws.Range("AQ3").Value = 685

# Row 4
$ws.Range("B4").Value = "African Rainbow Capital Investments Limited (JSE:AIL)"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.9465954606141521
$ws.Range("J4").Value = 0.9465954606141521
$ws.Range("K4").Value = 7.09
$ws.Range("L4").Value = 0.9465954606141521
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 1.1
$ws.Range("V4").Value = 0.002742458239840439
$ws.Range("W4").Value = 0.01013871013871014
$ws.Range("X4").Value = 0.0437464840940926
$ws.Range("Y4").Value = -0.03360777395538246
$ws.Range("Z4").Value = 0.01073681192660551
$ws.Range("AA4").Value = 0.01016341743119266
$ws.Range("AB4").Value = 0.0437464840940926
$ws.Range("AC4").Value = -0.03358306666289994
$ws.Range("AG4").Value = -1.1
$ws.Range("AJ4").Value = -0.00275
$ws.Range("AK4").Value = -0.001911713590545708
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("T4").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()

# Row 5
$ws.Range("B5").Value = "Bravura Holdings Limited (NMSE:CMB)"
$ws.Range("G5").Value = -0.1057507987220447
$ws.Range("H5").Value = -0.1057507987220447
$ws.Range("I5").Value = 0.004792332268370607
$ws.Range("J5").Value = 0.004792332268370607
$ws.Range("K5").Value = -19.8
$ws.Range("L5").Value = -3.162939297124601
$ws.Range("O5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("U5").Value = 0.773
$ws.Range("V5").Value = 0.0483125
$ws.Range("W5").Value = -0.5739130434782609
$ws.Range("X5").Value = 0.04449677157576705
$ws.Range("Y5").Value = -0.6184098150540279
$ws.Range("Z5").Value = 0.1849881796690307
$ws.Range("AA5").Value = 0.0008865248226950353
$ws.Range("AB5").Value = 0.04432140088924682
$ws.Range("AC5").Value = -0.04343487606655178
$ws.Range("AD5").Value = 0.41
$ws.Range("AF5").Value = 0.41
$ws.Range("AG5").Value = -0.363
$ws.Range("AH5").Value = 0.02498476538695917
$ws.Range("AI5").Value = 0.0299051787016776
$ws.Range("AJ5").Value = -0.0232141715162755
$ws.Range("AK5").Value = -0.02805905542243179
$ws.Range("AL5").Value = 0.04
$ws.Range("AM5").Value = -0.151
$ws.Range("AN5").Value = 0.9601873536299765
$ws.Range("AO5").Value = 0.75
$ws.Range("AP5").Value = -0.8501170960187355
$ws.Range("AQ5").Value = -0.1986754966887417

# Row 6
$ws.Range("B6").Value = "Promotion and Development Ltd (MUSE:PAD.N0000)"
$ws.Range("D6").Value = -0.00639
$ws.Range("G6").Value = 0.206015037593985
$ws.Range("H6").Value = 0.206015037593985
$ws.Range("I6").Value = 0.08721804511278194
$ws.Range("J6").Value = 0.08721804511278194
$ws.Range("K6").Value = -5.52
$ws.Range("L6").Value = -0.4150375939849624
$ws.Range("M6").Value = -0
$ws.Range("N6").Value = -0
$ws.Range("O6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("U6").Value = 0.08500000000000001
$ws.Range("V6").Value = 0.001125827814569536
$ws.Range("W6").Value = -0.01826000661594442
$ws.Range("X6").Value = 0.06146931422170884
$ws.Range("Y6").Value = -0.07972932083765327
$ws.Range("Z6").Value = 0.03792631458879891
$ws.Range("AA6").Value = 0.003307859016767423
$ws.Range("AB6").Value = 0.04896793313231951
$ws.Range("AC6").Value = -0.04566007411555209
$ws.Range("AD6").Value = 45.7
$ws.Range("AF6").Value = 45.7
$ws.Range("AG6").Value = 45.615
$ws.Range("AH6").Value = 0.3770627062706271
$ws.Range("AI6").Value = 0.1343327454438566
$ws.Range("AJ6").Value = 0.3766255211988606
$ws.Range("AK6").Value = 0.1341164018052718
$ws.Range("AL6").Value = 2.14
$ws.Range("AM6").Value = 2.14
$ws.Range("AO6").Value = 0.5420560747663551
$ws.Range("AQ6").Value = 0.5420560747663551
$ws.Range("T6").ClearContents()

# Row 7
$ws.Range("B7").Value = "The Bee Equity Partners Ltd (MUSE:FIDE.I0000)"
$ws.Range("D7").Value = -0.425
$ws.Range("G7").Value = 0.4147286821705426
$ws.Range("H7").Value = 0.4147286821705426
$ws.Range("I7").Value = -0.5155038759689923
$ws.Range("J7").Value = -0.5155038759689923
$ws.Range("K7").Value = -1.91
$ws.Range("L7").Value = -7.403100775193798
$ws.Range("M7").Value = 0.11466
$ws.Range("N7").Value = 0.02167485822306238
$ws.Range("O7").Value = -0.06003141361256545
$ws.Range("P7").Value = 0.11466
$ws.Range("Q7").Value = 0.02167485822306238
$ws.Range("R7").Value = -0.06003141361256545
$ws.Range("U7").Value = 2.8
$ws.Range("V7").Value = 0.5293005671077504
$ws.Range("W7").Value = -0.119375
$ws.Range("X7").Value = 0.0437464840940926
$ws.Range("Y7").Value = -0.1631214840940926
$ws.Range("Z7").Value = 0.01784232365145228
$ws.Range("AA7").Value = -0.009197786998616874
$ws.Range("AB7").Value = 0.0437464840940926
$ws.Range("AC7").Value = -0.05294427109270947
$ws.Range("AD7").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AG7").Value = -2.8
$ws.Range("AH7").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AJ7").Value = -1.124497991967871
$ws.Range("AK7").Value = -0.2772277227722772
$ws.Range("AL7").Value = 0.001
$ws.Range("AM7").Value = 0.001
$ws.Range("AN7").Value = -0
$ws.Range("AO7").Value = -133
$ws.Range("AP7").Value = 21.875
$ws.Range("AQ7").Value = -133
$ws.Range("E7").ClearContents()

# Row 8
$ws.Range("I8").Value = 1.102389078498293
$ws.Range("J8").Value = 1.102389078498293
$ws.Range("K8").Value = -3.22
$ws.Range("L8").Value = 1.098976109215017
$ws.Range("M8").Value = -0
$ws.Range("N8").Value = -0
$ws.Range("O8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("U8").Value = 0.619
$ws.Range("V8").Value = 0.01513447432762836
$ws.Range("W8").Value = -0.06110056925996205
$ws.Range("X8").Value = 0.0437464840940926
$ws.Range("Y8").Value = -0.1048470533540546
$ws.Range("Z8").Value = -0.05580527197927777
$ws.Range("AA8").Value = -0.06151912235258265
$ws.Range("AB8").Value = 0.0437464840940926
$ws.Range("AC8").Value = -0.1052656064466753
$ws.Range("AG8").Value = -0.619
$ws.Range("AJ8").Value = -0.0153670464983491
$ws.Range("AK8").Value = -0.012689366761649
$ws.Range("AM8").Value = -0.07099999999999999
$ws.Range("AQ8").Value = 45.49295774647888
$ws.Range("T8").ClearContents()

# Row 9
$ws.Range("B9").Value = "The Mauritius Development Investment Trust Company Limited (MUSE:MDIT.N0000)"
$ws.Range("G9").Value = -0
$ws.Range("H9").Value = -0
$ws.Range("I9").Value = 1.043329532497149
$ws.Range("J9").Value = 1.043329532497149
$ws.Range("K9").Value = -9.15
$ws.Range("L9").Value = 1.043329532497149
$ws.Range("M9").Value = 2.34
$ws.Range("N9").Value = 0.09322709163346612
$ws.Range("O9").Value = -0.2557377049180328
$ws.Range("P9").Value = 2.34
$ws.Range("Q9").Value = 0.09322709163346612
$ws.Range("R9").Value = -0.2557377049180328
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0.649
$ws.Range("V9").Value = 0.02585657370517928
$ws.Range("W9").Value = -0.2293233082706767
$ws.Range("X9").Value = 0.0437464840940926
$ws.Range("Y9").Value = -0.2730697923647693
$ws.Range("Z9").Value = -0.204739115209525
$ws.Range("AA9").Value = -0.2136103653554337
$ws.Range("AB9").Value = 0.0437464840940926
$ws.Range("AC9").Value = -0.2573568494495263
$ws.Range("AD9").Value = 0
$ws.Range("AE9").Value = 0
$ws.Range("AF9").Value = 0
$ws.Range("AG9").Value = -0.649
$ws.Range("AH9").Value = 0
$ws.Range("AI9").Value = 0
$ws.Range("AJ9").Value = -0.02654288168173081
$ws.Range("AK9").Value = -0.02491267129860658
$ws.Range("AL9").Value = 0.058
$ws.Range("AM9").Value = 0.058
$ws.Range("AO9").Value = -157.7586206896552
$ws.Range("AQ9").Value = -157.7586206896552
$ws.Range("AN9").ClearContents()
$ws.Range("AP9").ClearContents()

# Row 10
$ws.Range("B10").Value = "EPE Capital Partners Ltd (JSE:EPE)"
$ws.Range("I10").Value = 1.015555555555556
$ws.Range("J10").Value = 1.015555555555556
$ws.Range("K10").Value = -46.6
$ws.Range("L10").Value = 1.035555555555556
$ws.Range("U10").Value = 0.461
$ws.Range("V10").Value = 0.006557610241820768
$ws.Range("W10").Value = -0.3686708860759494
$ws.Range("X10").Value = 0.0482029673492245
$ws.Range("Y10").Value = -0.4168738534251739
$ws.Range("Z10").Value = -0.3315234608102433
$ws.Range("AA10").Value = -0.3366804924228471
$ws.Range("AB10").Value = 0.04557574739763805
$ws.Range("AC10").Value = -0.3822562398204851
$ws.Range("AD10").Value = 10.7
$ws.Range("AF10").Value = 10.7
$ws.Range("AG10").Value = 10.239
$ws.Range("AH10").Value = 0.1320987654320988
$ws.Range("AI10").Value = 0.09780621572212064
$ws.Range("AJ10").Value = 0.1271309551894113
$ws.Range("AK10").Value = 0.09398837881750335
$ws.Range("AL10").Value = 0.609
$ws.Range("AM10").Value = 0.609
$ws.Range("AO10").Value = -75.04105090311988
$ws.Range("AQ10").Value = -75.04105090311988
